# Auto-generated Excel COM-interop edit script
# Applies scraper refresh update (Linea 141) to sheets LP1912, LP1912-215, 6203-6173

$wb = $excel.ActiveWorkbook

# --- Sheet "LP1912": header + refreshed/re-sorted schedule rows ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = 'Última actualización: 08:46:25'

$ws1.Cells.Item(3,1).Value = 'Total filas: 111'

$ws1.Cells.Item(46,1).Value = '06:33:46'
$ws1.Cells.Item(46,2).Value = '07:59'
$ws1.Cells.Item(46,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(46,4).Value = 86
$ws1.Cells.Item(46,5).Value = 'LP1912'

$ws1.Cells.Item(47,1).Value = '07:12:53'
$ws1.Cells.Item(47,2).Value = '07:59'
$ws1.Cells.Item(47,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(47,4).Value = 47
$ws1.Cells.Item(47,5).Value = 'LP1912'

$ws1.Cells.Item(57,1).Value = '08:11:27'
$ws1.Cells.Item(57,2).Value = '08:14'
$ws1.Cells.Item(57,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(57,4).Value = 3
$ws1.Cells.Item(57,5).Value = 'LP1912'

$ws1.Cells.Item(58,1).Value = '07:48:35'
$ws1.Cells.Item(58,2).Value = '08:14'
$ws1.Cells.Item(58,3).Value = '10_OLMOS'
$ws1.Cells.Item(58,4).Value = 26
$ws1.Cells.Item(58,5).Value = 'LP1912'

$ws1.Cells.Item(59,1).Value = '07:36:59'
$ws1.Cells.Item(59,2).Value = '08:14'
$ws1.Cells.Item(59,3).Value = '17_ROMERO'
$ws1.Cells.Item(59,4).Value = 38
$ws1.Cells.Item(59,5).Value = 'LP1912'

$ws1.Cells.Item(97,1).Value = '08:46:25'
$ws1.Cells.Item(97,2).Value = '09:33'
$ws1.Cells.Item(97,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(97,4).Value = 47
$ws1.Cells.Item(97,5).Value = 'LP1912'

$ws1.Cells.Item(98,1).Value = '07:48:35'
$ws1.Cells.Item(98,2).Value = '09:39'
$ws1.Cells.Item(98,3).Value = '15_ABASTO'
$ws1.Cells.Item(98,4).Value = 111
$ws1.Cells.Item(98,5).Value = 'LP1912'

$ws1.Cells.Item(99,1).Value = '07:48:35'
$ws1.Cells.Item(99,2).Value = '09:44'
$ws1.Cells.Item(99,3).Value = '14_ABASTO'
$ws1.Cells.Item(99,4).Value = 116
$ws1.Cells.Item(99,5).Value = 'LP1912'

$ws1.Cells.Item(100,1).Value = '08:39:08'
$ws1.Cells.Item(100,2).Value = '09:45'
$ws1.Cells.Item(100,3).Value = '14_ABASTO'
$ws1.Cells.Item(100,4).Value = 66
$ws1.Cells.Item(100,5).Value = 'LP1912'

$ws1.Cells.Item(101,1).Value = '07:55:46'
$ws1.Cells.Item(101,2).Value = '09:51'
$ws1.Cells.Item(101,3).Value = '15_ABASTO'
$ws1.Cells.Item(101,4).Value = 116
$ws1.Cells.Item(101,5).Value = 'LP1912'

$ws1.Cells.Item(102,1).Value = '08:11:27'
$ws1.Cells.Item(102,2).Value = '10:03'
$ws1.Cells.Item(102,3).Value = '215C_EL PATO'
$ws1.Cells.Item(102,4).Value = 112
$ws1.Cells.Item(102,5).Value = 'LP1912'

$ws1.Cells.Item(103,1).Value = '08:46:25'
$ws1.Cells.Item(103,2).Value = '10:04'
$ws1.Cells.Item(103,3).Value = '14_ABASTO'
$ws1.Cells.Item(103,4).Value = 78
$ws1.Cells.Item(103,5).Value = 'LP1912'

$ws1.Cells.Item(104,1).Value = '08:39:08'
$ws1.Cells.Item(104,2).Value = '10:05'
$ws1.Cells.Item(104,3).Value = '14_ABASTO'
$ws1.Cells.Item(104,4).Value = 86
$ws1.Cells.Item(104,5).Value = 'LP1912'

$ws1.Cells.Item(105,1).Value = '08:11:27'
$ws1.Cells.Item(105,2).Value = '10:10'
$ws1.Cells.Item(105,3).Value = '10_OLMOS'
$ws1.Cells.Item(105,4).Value = 119
$ws1.Cells.Item(105,5).Value = 'LP1912'

$ws1.Cells.Item(106,1).Value = '08:29:19'
$ws1.Cells.Item(106,2).Value = '10:11'
$ws1.Cells.Item(106,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(106,4).Value = 102
$ws1.Cells.Item(106,5).Value = 'LP1912'

$ws1.Cells.Item(107,1).Value = '08:29:19'
$ws1.Cells.Item(107,2).Value = '10:12'
$ws1.Cells.Item(107,3).Value = '10_OLMOS'
$ws1.Cells.Item(107,4).Value = 103
$ws1.Cells.Item(107,5).Value = 'LP1912'

$ws1.Cells.Item(108,1).Value = '08:29:19'
$ws1.Cells.Item(108,2).Value = '10:14'
$ws1.Cells.Item(108,3).Value = '10_OLMOS'
$ws1.Cells.Item(108,4).Value = 105
$ws1.Cells.Item(108,5).Value = 'LP1912'

$ws1.Cells.Item(109,1).Value = '08:29:19'
$ws1.Cells.Item(109,2).Value = '10:15'
$ws1.Cells.Item(109,3).Value = '17_ROMERO'
$ws1.Cells.Item(109,4).Value = 106
$ws1.Cells.Item(109,5).Value = 'LP1912'

$ws1.Cells.Item(110,1).Value = '08:39:08'
$ws1.Cells.Item(110,2).Value = '10:15'
$ws1.Cells.Item(110,3).Value = '10_OLMOS'
$ws1.Cells.Item(110,4).Value = 96
$ws1.Cells.Item(110,5).Value = 'LP1912'

$ws1.Cells.Item(111,1).Value = '08:46:25'
$ws1.Cells.Item(111,2).Value = '10:16'
$ws1.Cells.Item(111,3).Value = '10_OLMOS'
$ws1.Cells.Item(111,4).Value = 90
$ws1.Cells.Item(111,5).Value = 'LP1912'

$ws1.Cells.Item(112,1).Value = '08:46:25'
$ws1.Cells.Item(112,2).Value = '10:18'
$ws1.Cells.Item(112,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(112,4).Value = 92
$ws1.Cells.Item(112,5).Value = 'LP1912'

$ws1.Cells.Item(113,1).Value = '08:29:19'
$ws1.Cells.Item(113,2).Value = '10:26'
$ws1.Cells.Item(113,3).Value = '15X38_ABASTO'
$ws1.Cells.Item(113,4).Value = 117
$ws1.Cells.Item(113,5).Value = 'LP1912'

$ws1.Cells.Item(114,1).Value = '08:39:08'
$ws1.Cells.Item(114,2).Value = '10:30'
$ws1.Cells.Item(114,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(114,4).Value = 111
$ws1.Cells.Item(114,5).Value = 'LP1912'

$ws1.Cells.Item(115,1).Value = '08:39:08'
$ws1.Cells.Item(115,2).Value = '10:34'
$ws1.Cells.Item(115,3).Value = '10_OLMOS'
$ws1.Cells.Item(115,4).Value = 115
$ws1.Cells.Item(115,5).Value = 'LP1912'

$ws1.Cells.Item(116,1).Value = '08:39:08'
$ws1.Cells.Item(116,2).Value = '10:37'
$ws1.Cells.Item(116,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(116,4).Value = 118
$ws1.Cells.Item(116,5).Value = 'LP1912'

# --- Sheet "LP1912-215": timestamp-only refresh ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = 'Última actualización: 08:46:25'

# --- Sheet "6203-6173": header + refreshed/re-sorted schedule rows ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = 'Última actualización: 08:46:25'

$ws3.Cells.Item(3,1).Value = 'Total filas: 14'

$ws3.Cells.Item(14,1).Value = '08:46:25'
$ws3.Cells.Item(14,2).Value = '08:48'
$ws3.Cells.Item(14,3).Value = '215A_LA PLATA'
$ws3.Cells.Item(14,4).Value = 2
$ws3.Cells.Item(14,5).Value = 'L6173'

$ws3.Cells.Item(15,1).Value = '07:36:59'
$ws3.Cells.Item(15,2).Value = '08:51'
$ws3.Cells.Item(15,3).Value = '215A_LA PLATA'
$ws3.Cells.Item(15,4).Value = 75
$ws3.Cells.Item(15,5).Value = 'L6173'

$ws3.Cells.Item(16,1).Value = '07:48:35'
$ws3.Cells.Item(16,2).Value = '08:52'
$ws3.Cells.Item(16,3).Value = '215A_LA PLATA'
$ws3.Cells.Item(16,4).Value = 64
$ws3.Cells.Item(16,5).Value = 'L6173'

$ws3.Cells.Item(17,1).Value = '08:11:27'
$ws3.Cells.Item(17,2).Value = '10:09'
$ws3.Cells.Item(17,3).Value = '215C_LA PLATA'
$ws3.Cells.Item(17,4).Value = 118
$ws3.Cells.Item(17,5).Value = 'L6203'

$ws3.Cells.Item(18,1).Value = '08:29:19'
$ws3.Cells.Item(18,2).Value = '10:23'
$ws3.Cells.Item(18,3).Value = '215A_LA PLATA'
$ws3.Cells.Item(18,4).Value = 114
$ws3.Cells.Item(18,5).Value = 'L6173'

$ws3.Cells.Item(19,1).Value = '08:39:08'
$ws3.Cells.Item(19,2).Value = '10:31'
$ws3.Cells.Item(19,3).Value = '215B_LP-P MOR-1 Y 57'
$ws3.Cells.Item(19,4).Value = 112
$ws3.Cells.Item(19,5).Value = 'L6173'

Write-Output "Updated LP1912, LP1912-215, 6203-6173 sheets."
